$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Helper: remove any hyperlink whose anchor range matches one of the given
# (dollar-sign) addresses, e.g. "$A$3". Leaves all other hyperlinks (and
# their XML ref="..." formatting) completely untouched.
#
# NOTE: the Hyperlinks collection re-indexes live as items are removed, so
# items must be collected first and then deleted back-to-front; deleting
# while walking forward causes entries to be skipped.
# ---------------------------------------------------------------------------
function Remove-HyperlinksAt($ws, $addresses) {
    $toDelete = @()
    foreach ($hl in $ws.Hyperlinks) {
        $addr = $hl.Range.Address()
        if ($addresses -contains $addr) {
            $toDelete += ,$hl
        }
    }
    for ($i = $toDelete.Count - 1; $i -ge 0; $i--) {
        $toDelete[$i].Delete()
    }
}

# ---------------------------------------------------------------------------
# Sheet 1 : "Overview"
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)

$ws1.Range("B2").Value = "Ready for handoff"
$ws1.Range("C2").Value = "Ready for handoff"
$ws1.Range("D2").Value = "2016-03-21 10:44:47"

Remove-HyperlinksAt $ws1 @("`$A`$3")
$ws1.Rows.Item(3).Delete()

# ---------------------------------------------------------------------------
# Sheet 2 : "zh-cn"
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item(2)

$ws2.Range("C2").Value = "Ready for handoff"
$ws2.Range("E2").Value = "2016-03-21 10:44:43"

Remove-HyperlinksAt $ws2 @("`$A`$3", "`$D`$3", "`$F`$3", "`$G`$3")
$ws2.Rows.Item(3).Delete()

# ---------------------------------------------------------------------------
# Sheet 3 : "de-de"
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item(3)

$ws3.Range("C2").Value = "Ready for handoff"
$ws3.Range("E2").Value = "2016-03-21 10:44:47"

Remove-HyperlinksAt $ws3 @("`$A`$3", "`$D`$3", "`$F`$3", "`$G`$3")
$ws3.Rows.Item(3).Delete()
